$d = $word.ActiveDocument

# Update the date line at the top of the document. This string is unique in
# the document, so a plain Find/Replace is safe here.
$d.Content.Find.Execute("2026-02-13 Friday", $true, $false, $false, $false,
    $false, $true, 1, $false, "2026-02-14 Saturday", 2)

# Update each arithmetic-problem cell in the table, addressed by explicit
# (row, column) position rather than Find/Replace. Several of the new
# values collide with old values that live elsewhere in the table (e.g.
# "17÷7=" is being replaced in one cell but introduced as a new value in
# another, and likewise for "68÷6=" / "24÷6=" / "38÷9="), so a blanket
# text search-and-replace would cause cascading/incorrect substitutions.
# Targeting cells directly avoids that entirely.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92÷5="
$t.Cell(1, 2).Range.Text = "25÷4="
$t.Cell(1, 3).Range.Text = "38÷9="
$t.Cell(1, 4).Range.Text = "41÷6="
$t.Cell(1, 5).Range.Text = "84÷8="

$t.Cell(5, 1).Range.Text = "92÷3="
$t.Cell(5, 2).Range.Text = "64÷9="
$t.Cell(5, 3).Range.Text = "44÷3="
$t.Cell(5, 4).Range.Text = "24÷6="
$t.Cell(5, 5).Range.Text = "32÷2="

$t.Cell(9, 1).Range.Text = "35÷6="
$t.Cell(9, 2).Range.Text = "39÷5="
$t.Cell(9, 3).Range.Text = "45÷2="
$t.Cell(9, 4).Range.Text = "62÷8="
$t.Cell(9, 5).Range.Text = "10÷9="

$t.Cell(13, 1).Range.Text = "46÷3="
$t.Cell(13, 2).Range.Text = "95÷8="
$t.Cell(13, 3).Range.Text = "78÷5="
$t.Cell(13, 4).Range.Text = "38÷9="
$t.Cell(13, 5).Range.Text = "17÷7="

$t.Cell(17, 1).Range.Text = "37÷8="
$t.Cell(17, 2).Range.Text = "24÷6="
$t.Cell(17, 3).Range.Text = "51÷9="
$t.Cell(17, 4).Range.Text = "68÷6="
$t.Cell(17, 5).Range.Text = "55÷6="
